# Update the LR-pairs sheet with the new TPM-derived values.
# The data table (Sending cluster x Target cluster combinations for the
# Bdnf -> Ntrk2 ligand/receptor pair) grows from 3 rows to 6 rows, covering
# every combination of Sending cluster in {ECs, MuSCs} x Target cluster in
# {ECs, FAPs, MuSCs}, and every numeric column is refreshed with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  A="ECs";   B="Bdnf"; C="Ntrk2"; D="ECs";
       E=1;  F=0.3333333333333333; G=0.03885866666666667; H=0.116576;
       I=0.09340097618505853; J=0.09340097618505853;
       K=3;  L=1; M=0.165174; N=0.495522;
       O=0.00279520163427027; P=0.002795201634270271;
       Q=0.006418441408000001; R=0.057765972672;
       S=0.0002610745612749142; T=0.0002610745612749142 },

    @{ Row=3;  A="ECs";   B="Bdnf"; C="Ntrk2"; D="FAPs";
       E=1;  F=0.3333333333333333; G=0.03885866666666667; H=0.116576;
       I=0.09340097618505853; J=0.09340097618505853;
       K=3;  L=1; M=45.53127133333334; N=136.593814;
       O=0.7705152387260491; P=0.7705152387260492;
       Q=1.769284495651556; R=15.923560460864;
       S=0.0719668754624764; T=0.07196687546247642 },

    @{ Row=4;  A="ECs";   B="Bdnf"; C="Ntrk2"; D="MuSCs";
       E=1;  F=0.3333333333333333; G=0.03885866666666667; H=0.116576;
       I=0.09340097618505853; J=0.09340097618505853;
       K=3;  L=1; M=13.395535; N=40.186605;
       O=0.2266895596396806; P=0.2266895596396806;
       Q=0.5205326293866667; R=4.68479366448;
       S=0.02117302616130722; T=0.02117302616130722 },

    @{ Row=5;  A="MuSCs"; B="Bdnf"; C="Ntrk2"; D="ECs";
       E=2;  F=0.6666666666666666; G=0.3771826666666667; H=1.131548;
       I=0.9065990238149415; J=0.9065990238149415;
       K=3;  L=1; M=0.165174; N=0.495522;
       O=0.00279520163427027; P=0.002795201634270271;
       Q=0.062300769784; R=0.560706928056;
       S=0.002534127072995356; T=0.002534127072995356 },

    @{ Row=6;  A="MuSCs"; B="Bdnf"; C="Ntrk2"; D="FAPs";
       E=2;  F=0.6666666666666666; G=0.3771826666666667; H=1.131548;
       I=0.9065990238149415; J=0.9065990238149415;
       K=3;  L=1; M=45.53127133333334; N=136.593814;
       O=0.7705152387260491; P=0.7705152387260492;
       Q=17.17360633823022; R=154.562457044072;
       S=0.6985483632635727; T=0.6985483632635728 },

    @{ Row=7;  A="MuSCs"; B="Bdnf"; C="Ntrk2"; D="MuSCs";
       E=2;  F=0.6666666666666666; G=0.3771826666666667; H=1.131548;
       I=0.9065990238149415; J=0.9065990238149415;
       K=3;  L=1; M=13.395535; N=40.186605;
       O=0.2266895596396806; P=0.2266895596396806;
       Q=5.052563612726667; R=45.47307251454;
       S=0.2055165334783734; T=0.2055165334783734 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
